$d = $word.ActiveDocument
$sel = $word.Selection
$sel.EndKey(6)  # wdStory: move the insertion point to the very end of the document

function Insert-TrackedBreak {
    param($selection, $doc)
    $before = $doc.Content.End
    $selection.InsertBreak(6)  # wdLineBreak -> <w:br/>
    $after = $doc.Content.End
    $rng = $doc.Range($before - 1, $after - 1)
    $rng.Font.Size = 14  # 14pt == w:sz val="28" (half-points); forces explicit rPr on the new run
}

function Insert-TrackedText {
    param($selection, $doc, $text)
    $before = $doc.Content.End
    $selection.TypeText($text)
    $after = $doc.Content.End
    $rng = $doc.Range($before - 1, $after - 1)
    $rng.Font.Size = 14
}

Insert-TrackedBreak $sel $d
Insert-TrackedBreak $sel $d
Insert-TrackedText $sel $d "https://share.weiyun.com/Gh1Kdv4z"
Insert-TrackedBreak $sel $d
